$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Techno (B), Personnel (C), Travaux (D) columns for rows 3-15
# reflecting the new DB structure / foreign key / project plan changes.

# Row 3
$ws.Range("B3").Value = "HTML + PHP"
$ws.Range("C3").Value = "Moi"
$ws.Range("D3").Value = "Création de structure MVC + configuration BDD"

# Row 4
$ws.Range("B4").Value = "HTML + PHP + MySQL"
$ws.Range("C4").Value = "Moi"
$ws.Range("D4").Value = "Connexion Backend + sécurisation & gestions des erreurs"

# Row 5
$ws.Range("B5").Value = "HTML + PHP + MySQL"
$ws.Range("C5").Value = "Moi"
$ws.Range("D5").Value = "CRUD Admin"

# Row 6
$ws.Range("B6").Value = "HTML + PHP + MySQL"
$ws.Range("C6").Value = "Moi"
$ws.Range("D6").Value = "CRUD Admin"

# Row 7
$ws.Range("B7").Value = "HTML + PHP + MySQL"
$ws.Range("C7").Value = "Moi"
$ws.Range("D7").Value = "CRUD Content"

# Row 8
$ws.Range("B8").Value = "HTML + JS + CSS"
$ws.Range("C8").Value = "Moi"
$ws.Range("D8").Value = "FRONT Barre nav, structure de page"

# Row 9
$ws.Range("B9").Value = "HTML + JS + CSS"
$ws.Range("C9").Value = "Moi"
$ws.Range("D9").Value = "FRONT Footer, Charte graphique"

# Row 10
$ws.Range("B10").Value = "HTML + JS + CSS + PHP"
$ws.Range("C10").Value = "Moi"
$ws.Range("D10").Value = "FRONT Page Salon  (Dynamique)"

# Row 11
$ws.Range("B11").Value = "HTML + JS + CSS + PHP"
$ws.Range("C11").Value = "Moi"
$ws.Range("D11").Value = "FRONT Page Salon  (Dynamique)"

# Row 12
$ws.Range("B12").Value = "HTML + JS + CSS"
$ws.Range("C12").Value = "Moi"
$ws.Range("D12").Value = "FRONT Page Equipe (Statique)"

# Row 13
$ws.Range("B13").Value = "HTML + JS + CSS + PHP"
$ws.Range("C13").Value = "Moi"
$ws.Range("D13").Value = "FRONT Page Actualité (Dynamique)"

# Row 14
$ws.Range("B14").Value = "HTML + JS + CSS + PHP"
$ws.Range("C14").Value = "Moi"
$ws.Range("D14").Value = "FRONT Page Actualité (Dynamique)"

# Row 15
$ws.Range("B15").Value = "HTML + JS + CSS + PHP"
$ws.Range("C15").Value = "Moi"
$ws.Range("D15").Value = "FRONT Page Actualité (Dynamique)"

# Update the active selection to match the final cursor position
$ws.Range("C23").Select()
